$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.991.36"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.919.57"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'325.58"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4590"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.3812"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "'0.07751"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "'0.9772"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'22.63"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("D12").Value = "1.934.12"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "'5.708"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'6.945"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "'0.07028"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'84.55"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'0.000009500"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "29.015.02"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'5.347"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "2.180.07"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "'157.78"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "'5.608"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "'117.75"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "'1.831"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'0.09327"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "'0.8593"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'5.088"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "'3.024"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'0.05679"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "'1.151"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'1.004"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'3.142"
$ws.Range("E39").Value = "  +14.69%  "
$ws.Range("D40").Value = "'0.02041"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "'7.427"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").Value = "'0.5489"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'9.363"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "'0.000002844"
$ws.Range("E45").Value = "  +8.63%  "
$ws.Range("D46").Value = "'2.178"
$ws.Range("E46").Value = "  +3.86%  "
$ws.Range("D47").Value = "'0.5182"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "'11.21"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").Value = "'0.06915"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").Value = "'110.22"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("E51").Value = "  -1.23%  "
